$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns remain plain text so values such as
# "1.00", "304.81", "42.755.95" are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.755.95"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.523.09"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "304.81"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "96.76"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "0.579"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").Value = "36.44"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "0.0811"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "7.69"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "2.907.49"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "2.511.18"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "15.09"
$ws.Range("E16").Value = "  +4.88%  "
$ws.Range("D17").Value = "0.860"
$ws.Range("E17").Value = "  -3.12%  "
$ws.Range("D18").Value = "42.705.74"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "12.98"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "0.0₃0975"
$ws.Range("E20").Value = "  -2.27%  "
$ws.Range("D21").Value = "6.46"
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("D22").Value = "71.21"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -3.96%  "
$ws.Range("D26").Value = "27.07"
$ws.Range("E26").Value = "  -6.44%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "2.32"
$ws.Range("E28").Value = "  +10.26%  "
$ws.Range("D29").Value = "10.32"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "38.27"
$ws.Range("E30").Value = "  +1.51%  "
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "155.97"
$ws.Range("E32").Value = "  +0.61%  "
$ws.Range("E33").Value = "  -3.35%  "
$ws.Range("D34").Value = "0.0790"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  -4.86%  "
$ws.Range("D36").Value = "2.63"
$ws.Range("E36").Value = "  -4.49%  "
$ws.Range("D37").Value = "18.55"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "24.16"
$ws.Range("E39").Value = "  +4.29%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "0.119"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "2.06"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "3.85"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("D46").Value = "2.035.03"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("D47").Value = "85.38"
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("D48").Value = "8.94"
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("D49").Value = "2.767.67"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "0.190"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").Value = "101.95"
$ws.Range("E51").Value = "  -4.26%  "
